$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 in the "Sample Project" rules sheet currently holds the text "R40"
# in B11. The edit replaces that cell's text with "1" (still a text label,
# not a number, matching the other rule-name cells in column B such as
# "R10"/"R20"/"R30"). A leading apostrophe forces Excel to store it as text
# (shared string) rather than auto-converting it to a numeric value.
$ws.Range("B11").Value = "'1"
